$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.712.37"
$ws.Range("E2").Value = "  +0.36%  "
$ws.Range("D3").Value = "1.599.58"
$ws.Range("E3").Value = "  +0.20%  "
$ws.Range("E4").Value = "  +0.30%  "
$ws.Range("D5").Value = "'211.42"
$ws.Range("E5").Value = "  -0.08%  "
$ws.Range("E6").Value = "  -0.68%  "
$ws.Range("E7").Value = "  +0.29%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("E9").Value = "  +0.79%  "
$ws.Range("D10").Value = "'19.53"
$ws.Range("E10").Value = "  +0.26%  "
$ws.Range("E11").Value = "  +0.54%  "
$ws.Range("D12").Value = "1.824.82"
$ws.Range("E12").Value = "  +0.25%  "
$ws.Range("D13").Value = "1.623.08"
$ws.Range("E13").Value = "  +2.11%  "
$ws.Range("E14").Value = "  +0.51%  "
$ws.Range("D15").Value = "'0.522"
$ws.Range("E15").Value = "  +0.14%  "
$ws.Range("D16").Value = "'65.27"
$ws.Range("E16").Value = "  +1.28%  "
$ws.Range("D17").Value = "26.688.05"
$ws.Range("E17").Value = "  +0.30%  "
$ws.Range("D18").Value = "0.0₃0754"
$ws.Range("E18").Value = "  +3.08%  "
$ws.Range("D19").Value = "'7.20"
$ws.Range("E19").Value = "  +3.51%  "
$ws.Range("E20").Value = "  +0.33%  "
$ws.Range("D21").Value = "'208.98"
$ws.Range("E21").Value = "  +0.25%  "
$ws.Range("E22").Value = "  +0.46%  "
$ws.Range("E23").Value = "  +0.09%  "
$ws.Range("D24").Value = "'8.93"
$ws.Range("E24").Value = "  +0.52%  "
$ws.Range("D25").Value = "'142.27"
$ws.Range("E25").Value = "  -2.03%  "
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("D27").Value = "'7.11"
$ws.Range("E27").Value = "  -0.29%  "
$ws.Range("D28").Value = "'0.114"
$ws.Range("E28").Value = "  -0.04%  "
$ws.Range("D29").Value = "'15.34"
$ws.Range("E29").Value = "  +0.50%  "
$ws.Range("E30").Value = "  +2.66%  "
$ws.Range("E31").Value = "  -0.37%  "
$ws.Range("E32").Value = "  +0.67%  "
$ws.Range("E33").Value = "  +1.27%  "
$ws.Range("D34").Value = "1.291.03"
$ws.Range("E34").Value = "  +0.75%  "
$ws.Range("E35").Value = "  -5.34%  "
$ws.Range("E36").Value = "  +0.96%  "
$ws.Range("E37").Value = "  +0.07%  "
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("E39").Value = "  +20.30%  "
$ws.Range("E40").Value = "  -1.98%  "
$ws.Range("D41").Value = "'5.41"
$ws.Range("E41").Value = "  -1.10%  "
$ws.Range("E42").Value = "  +0.06%  "
$ws.Range("D43").Value = "'0.784"
$ws.Range("E43").Value = "  -0.21%  "
$ws.Range("D44").Value = "'63.13"
$ws.Range("E44").Value = "  -1.97%  "
$ws.Range("D45").Value = "1.736.86"
$ws.Range("E45").Value = "  +0.21%  "
$ws.Range("D46").Value = "'91.22"
$ws.Range("D47").Value = "'1.57"
$ws.Range("E47").Value = "  -1.68%  "
$ws.Range("E48").Value = "  -2.18%  "
$ws.Range("E49").Value = "  +0.66%  "
$ws.Range("E50").Value = "  +0.36%  "
$ws.Range("D51").Value = "'7.39"
$ws.Range("E51").Value = "  -1.31%  "
